# Fill column E (rows 1-92) on Sheet4 with "x" for every row (mirrors the
# already-filled column D), and leave the final selection on E90 to match
# the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Activate() | Out-Null

$ws.Range("E1:E92").Value = "x"

$ws.Range("E90").Select() | Out-Null
